$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the Price/Volume columns so numeric-looking
# strings (e.g. '1.000', '25.763.68') are not coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$updates = [ordered]@{
    'D2' = '25.763.68'
    'E2' = '  -4.50%  '
    'D3' = '1.813.59'
    'E3' = '  -3.44%  '
    'D4' = '1.000'
    'E4' = '  -0.28%  '
    'D5' = '276.58'
    'E5' = '  -8.44%  '
    'D6' = '0.9998'
    'E6' = '  -0.25%  '
    'D7' = '0.5089'
    'E7' = '  -4.04%  '
    'D8' = '0.3509'
    'E8' = '  -6.67%  '
    'D9' = '44.80'
    'E9' = '  -1.67%  '
    'D10' = '0.06670'
    'E10' = '  -6.93%  '
    'D11' = '19.95'
    'E11' = '  -7.76%  '
    'D12' = '0.8298'
    'E12' = '  -6.37%  '
    'D13' = '0.07872'
    'E13' = '  -3.08%  '
    'D14' = '1.802.48'
    'E14' = '  -3.76%  '
    'D15' = '5.065'
    'E15' = '  -4.11%  '
    'D16' = '87.57'
    'E16' = '  -5.65%  '
    'D17' = '0.9992'
    'E17' = '  -0.40%  '
    'D18' = '13.95'
    'E18' = '  -5.49%  '
    'D19' = '0.000008056'
    'E19' = '  -5.96%  '
    'D20' = '0.9999'
    'E20' = '  -0.09%  '
    'D21' = '25.795.71'
    'E21' = '  -4.96%  '
    'D22' = '4.716'
    'E22' = '  -4.99%  '
    'D23' = '9.972'
    'E23' = '  -7.05%  '
    'D24' = '6.038'
    'E24' = '  -5.78%  '
    'E25' = '  -2.77%  '
    'D26' = '140.16'
    'E26' = '  -4.91%  '
    'D27' = '1.663'
    'E27' = '  -4.42%  '
    'D28' = '17.01'
    'E28' = '  -5.74%  '
    'D29' = '109.37'
    'E29' = '  -4.62%  '
    'D30' = '4.325'
    'E30' = '  -8.50%  '
    'D31' = '4.226'
    'E31' = '  -7.63%  '
    'D32' = '0.08786'
    'E32' = '  -3.51%  '
    'D33' = '0.04863'
    'E33' = '  -2.41%  '
    'D34' = '1.143'
    'E34' = '  -2.47%  '
    'B35' = 'ImmutableX'
    'C35' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'D35' = '0.7254'
    'E35' = '  -8.98%  '
    'B36' = 'HuobiToken'
    'C36' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    'D36' = '2.873'
    'E36' = '  -4.90%  '
    'D37' = '0.9993'
    'E37' = '  -0.19%  '
    'D38' = '3.134'
    'E38' = '  -3.02%  '
    'D39' = '0.5186'
    'E39' = '  -11.37%  '
    'D40' = '0.01839'
    'E40' = '  -5.98%  '
    'D41' = '2.245'
    'E41' = '  -13.75%  '
    'D42' = '0.9504'
    'E42' = '  -11.07%  '
    'D43' = '113.03'
    'E43' = '  -2.62%  '
    'D44' = '6.136'
    'E44' = '  -6.65%  '
    'D45' = '8.047'
    'E45' = '  -9.38%  '
    'D46' = '0.9993'
    'E46' = '  -0.19%  '
    'B47' = 'Decentraland'
    'C47' = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
    'D47' = '0.4537'
    'E47' = '  -9.84%  '
    'B48' = 'Algorand'
    'C48' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'D48' = '0.1360'
    'E48' = '  -9.18%  '
    'D49' = '9.266'
    'E49' = '  -6.67%  '
    'D50' = '36.17'
    'E50' = '  -4.09%  '
    'D51' = '1.498'
    'E51' = '  -6.63%  '
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

"Updated " + $updates.Count + " cells."
